$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: remove the stray _GoBack bookmark that sits after "... entero" in
# the resource-description paragraph near the top of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Edit 2: shorten the long instruction paragraph - the "En cada enunciado..."
# sentence moves out of this paragraph (it becomes its own paragraph, see
# Edit 3 below).
# ---------------------------------------------------------------------------
$oldInstruction = "Construye, en tu cuaderno, una recta numérica desde ‒7 hasta 9. En cada enunciado identifica el número entero que está inmediatamente a la derecha o izquierda del número que allí se indica. Luego, arrastra el número de la columna derecha frente a la frase que completa cada enunciado."
$newInstruction = "Construye, en tu cuaderno, una recta numérica desde ‒7 hasta 9. Luego, arrastra el número de la columna derecha frente a la frase que completa cada enunciado."

$d.Content.Find.Execute($oldInstruction, $true, $true, $false, $false, $false, $true, 1, $false, $newInstruction, 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: the placeholder paragraph that used to just contain "N" (right
# after "Más información (ventana flotante)") now carries the sentence that
# was removed from the instruction above, split across three runs:
#   1) "En cada enunciado ... se indica"
#   2) ". "
#   3) "Luego, arrastra ... enunciado."
# All three runs keep the same run formatting as the original "N" run.
# ---------------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr>'

$part1 = "En cada enunciado identifica el número entero que está inmediatamente a la derecha o izquierda del número que allí se indica"
$part2 = ". "
$part3 = "Luego, arrastra el número de la columna derecha frente a la frase que completa cada enunciado."

# Locate the target paragraph: it is the one right after the paragraph that
# contains "(ventana flotante)" and whose own text is exactly "N".
$targetPara = $null
$prevWasFlotante = $false
foreach ($p in $d.Paragraphs) {
    if ($prevWasFlotante -and $p.Range.Text -eq "N`r") {
        $targetPara = $p
        break
    }
    $prevWasFlotante = $p.Range.Text -like "*(ventana flotante)*"
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""

    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="399E0018" w14:textId="6CE812EB" w:rsidR="00D073C3" w:rsidRPr="000719EE" w:rsidRDefault="00B0726F" w:rsidP="00D073C3"><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>$part1</w:t></w:r><w:r>$rPr<w:t xml:space="preserve">$part2</w:t></w:r><w:r>$rPr<w:t>$part3</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Edit 4: put the _GoBack bookmark back, this time right after the closing
# ")" run of the "Sin ordenación aleatoria (S/N):)" paragraph (collapsed,
# i.e. bookmarkStart immediately followed by bookmarkEnd).
#
# A directly-collapsed Range positioned exactly on a paragraph-mark boundary
# confuses Bookmarks.Add in this host, so we insert a throwaway marker run,
# bookmark the (now safely mid-paragraph) range, and delete the marker text
# again - leaving the bookmark collapsed in the right spot.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Sin ordenación aleatoria (S/N):)", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("ZZTEMPMARKERZZ")
    $d.Bookmarks.Add("_GoBack", $rng)
    $rng.Text = ""
}
